$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) First paragraph: "This is a Microsoft word document." gets two trailing
#    spaces appended, followed by a red parenthetical note split across
#    three runs: "(This is a change \u2013 Ve" | "rsion for main branch" | ")"
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)

# Append two trailing spaces to the existing sentence (keeps it in the same run).
$end = $p1.Range.End - 1
$r = $d.Range($end, $end)
$r.InsertAfter("  ")

# Run 2 (red): "(This is a change – Ve"
$insPos = $p1.Range.End - 1
$rA = $d.Range($insPos, $insPos)
$rA.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$afterA = $p1.Range.End - 1
$colorA = $d.Range($insPos, $afterA)
$colorA.Font.Color = 255

# Run 3 (red): "rsion for main branch"
$insPos2 = $p1.Range.End - 1
$rB = $d.Range($insPos2, $insPos2)
$rB.InsertAfter("rsion for main branch")
$afterB = $p1.Range.End - 1
$colorB = $d.Range($insPos2, $afterB)
$colorB.Font.Color = 255

# Run 4 (red): ")"
$insPos3 = $p1.Range.End - 1
$rC = $d.Range($insPos3, $insPos3)
$rC.InsertAfter(")")
$afterC = $p1.Range.End - 1
$colorC = $d.Range($insPos3, $afterC)
$colorC.Font.Color = 255

# ---------------------------------------------------------------------------
# 2) Append a new, bare, shaded paragraph after the document's final
#    paragraph ("... we are free at last.").
# ---------------------------------------------------------------------------
$endOfDoc = $d.Content.End
$tailRange = $d.Range($endOfDoc, $endOfDoc)

$shadedParaPackage = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$tailRange.InsertXML($shadedParaPackage)

Write-Output "done"
